$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sales data (rows 2-15) replacing the previous 31/1-3/2 entries; row 16 is dropped
$data = @(
    @('V-1770508238446', '7/2/2026', '06:50 p. m.', 'Laura', 'Cerveza Corona (x1)', 10000, 1),
    @('V-1770508277150', '7/2/2026', '06:51 p. m.', 'Laura', 'Cerveza Corona (x1)', 10000, 2),
    @('V-1770509831306', '7/2/2026', '07:17 p. m.', 'Laura', 'Aguardiente Amarillo Media (x1)', 70000, 0),
    @('V-1770510307354', '7/2/2026', '07:25 p. m.', 'Laura', 'Cerveza Corona (x8)', 80000, 0),
    @('V-1770510578344', '7/2/2026', '07:29 p. m.', 'Laura', 'Cerveza Corona (x1), Ron viejo de caldas (5años) botella (x1), Aguardiente Amarillo Media (x1)', 212000, 2),
    @('V-1770511453742', '7/2/2026', '07:44 p. m.', 'Laura', 'Cerveza Corona (x1)', 10000, 2),
    @('V-1770511491261', '7/2/2026', '07:44 p. m.', 'Laura', 'Cerveza Corona (x1)', 10000, 0),
    @('V-1770511516090', '7/2/2026', '07:45 p. m.', 'Laura', 'Cerveza Corona (x6)', 60000, 0),
    @('V-1770511548532', '7/2/2026', '07:45 p. m.', 'Laura', 'Cerveza Corona (x3)', 30000, 0),
    @('V-1770512051170', '7/2/2026', '07:54 p. m.', 'Laura', 'Cerveza Corona (x1)', 10000, 0),
    @('V-1770512294788', '7/2/2026', '07:58 p. m.', 'Laura', 'Ron viejo de caldas (5años) botella (x1), Aguardiente Amarillo Media (x1), Cerveza Corona (x1)', 212000, 2),
    @('V-1770513802493', '7/2/2026', '08:23 p. m.', 'Laura', 'Cerveza Corona (x1), Ron viejo de caldas (5años) botella (x1)', 142000, 2),
    @('V-1770519457303', '7/2/2026', '09:57 p. m.', 'Laura', 'Cerveza Corona (x1), Aguardiente Amarillo Media (x1), Ron viejo de caldas (5años) botella (x1)', 212000, 1),
    @('V-1770519511759', '7/2/2026', '09:58 p. m.', 'Laura', 'Cerveza Corona (x1), Ron viejo de caldas (5años) botella (x1), Aguardiente Amarillo Media (x1)', 212000, 0)
)

# Keep the Fecha column (B) as plain text like the rest of the sheet,
# rather than letting Excel coerce "7/2/2026" into a date serial: mark the
# cells as text before writing, then restore the default "Normal" style so
# no stray number-format style lingers on the cell.
$ws.Range("B2:B15").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("B2:B15").Style = "Normal"

# Remove the now-obsolete row 16 (data shrank from 16 to 15 rows)
$ws.Rows.Item(16).Delete()

Write-Host "Done. New dimension: $($ws.UsedRange.Address())"
